$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the R (C2) and X (D2) values were entered in the wrong columns -
# swap them so C2/D2 hold the correct R/X figures.
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$ws.Range("C2").Value2 = $d2
$ws.Range("D2").Value2 = $c2

# Imax column (E2:E34): update line rating from 2000 to 1000.
$ws.Range("E2:E34").Value2 = 1000
